$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9660525321960449
$ws.Range("B1").Value = 1.904026985168457
$ws.Range("C1").Value = 4.966778755187988
$ws.Range("D1").Value = 2.010943651199341
$ws.Range("E1").Value = 0.6163186430931091
